$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; this shifts the existing D:K data (and
# formatting) one column to the right, into E:L, leaving a blank column D.
$ws.Columns("D").Insert()

# The freshly inserted column D picked up column C's formatting (Excel's
# "insert" default). Re-apply the per-row number formatting that column D
# should actually have -- i.e. the same formatting now sitting in column E
# (which is exactly what column D had before the insert).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new (most recent) period's figures in column D.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 135800
$ws.Range("D15").Value = -1000
$ws.Range("D17").Value = 9300
$ws.Range("D18").Value = 126400
$ws.Range("D20").Value = -67800
$ws.Range("D21").Value = 62900
$ws.Range("D23").Value = 58700
$ws.Range("D24").Value = 11400
$ws.Range("D26").Value = 47300
$ws.Range("D27").Value = 47300
$ws.Range("D29").Value = 0
$ws.Range("D32").Value = 67800
$ws.Range("D33").Value = 47300
$ws.Range("D35").Value = 47300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 54500
$ws.Range("D42").Value = 340600
$ws.Range("D48").Value = 51600
$ws.Range("D49").Value = 86700
$ws.Range("D52").Value = 7200
$ws.Range("D54").Value = 3279100
$ws.Range("D61").Value = 1600
$ws.Range("D66").Value = 2791500
$ws.Range("D72").Value = 160600
$ws.Range("D76").Value = 487600
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 47300
$ws.Range("D83").Value = 4300
$ws.Range("D89").Value = 49300
$ws.Range("D91").Value = -1300
$ws.Range("D94").Value = -147700
$ws.Range("D96").Value = -5000
$ws.Range("D100").Value = 154300
$ws.Range("D102").Value = 55900
